$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'45.726.58"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -1.27%  '
$ws.Cells.Item(3, 4).Value = "'2.592.41"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -0.71%  '
$ws.Cells.Item(4, 5).Value = '  +0.10%  '
$ws.Cells.Item(5, 4).Value = "'308.29"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -1.48%  '
$ws.Cells.Item(6, 4).Value = "'98.53"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -2.11%  '
$ws.Cells.Item(7, 5).Value = '  -0.96%  '
$ws.Cells.Item(8, 5).Value = '  +0.12%  '
$ws.Cells.Item(9, 4).Value = "'0.574"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -1.31%  '
$ws.Cells.Item(10, 4).Value = "'38.47"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -0.75%  '
$ws.Cells.Item(11, 4).Value = "'53.97"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.00%  '
$ws.Cells.Item(12, 4).Value = "'0.0836"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -0.71%  '
$ws.Cells.Item(13, 4).Value = "'8.02"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -3.87%  '
$ws.Cells.Item(14, 4).Value = "'2.993.91"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -0.65%  '
$ws.Cells.Item(15, 5).Value = '  +0.67%  '
$ws.Cells.Item(16, 4).Value = "'2.605.36"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.21%  '
$ws.Cells.Item(17, 4).Value = "'0.907"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.24%  '
$ws.Cells.Item(18, 4).Value = "'14.69"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -1.99%  '
$ws.Cells.Item(19, 4).Value = "'45.807.17"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -1.43%  '
$ws.Cells.Item(20, 5).Value = '  -1.00%  '
$ws.Cells.Item(21, 4).Value = "'6.66"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -1.46%  '
$ws.Cells.Item(22, 4).Value = "'12.58"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -5.14%  '
$ws.Cells.Item(23, 4).Value = "'284.74"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +11.42%  '
$ws.Cells.Item(24, 4).Value = "'73.18"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +2.87%  '
$ws.Cells.Item(25, 5).Value = '  -2.40%  '
$ws.Cells.Item(26, 5).Value = '  -0.02%  '
$ws.Cells.Item(27, 4).Value = "'28.99"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +2.48%  '
$ws.Cells.Item(28, 5).Value = '  +0.24%  '
$ws.Cells.Item(29, 5).Value = '  +0.48%  '
$ws.Cells.Item(30, 4).Value = "'10.58"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -0.29%  '
$ws.Cells.Item(31, 4).Value = "'38.28"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -5.91%  '
$ws.Cells.Item(32, 5).Value = '  -3.38%  '
$ws.Cells.Item(33, 5).Value = '  +0.55%  '
$ws.Cells.Item(34, 4).Value = "'3.59"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -3.42%  '
$ws.Cells.Item(35, 4).Value = "'157.79"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +3.11%  '
$ws.Cells.Item(36, 5).Value = '  -2.43%  '
$ws.Cells.Item(37, 4).Value = "'2.80"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -3.02%  '
$ws.Cells.Item(38, 5).Value = '  -1.27%  '
$ws.Cells.Item(39, 4).Value = "'0.121"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +2.22%  '
$ws.Cells.Item(40, 4).Value = "'0.123"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.06%  '
$ws.Cells.Item(41, 4).Value = "'15.74"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -7.44%  '
$ws.Cells.Item(42, 4).Value = "'0.0325"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.62%  '
$ws.Cells.Item(43, 2).Value = 'EnergySwap'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(43, 4).Value = "'21.40"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +1.49%  '
$ws.Cells.Item(44, 2).Value = 'NEARProtocol'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(44, 4).Value = "'3.51"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -2.83%  '
$ws.Cells.Item(45, 2).Value = 'RenderToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(45, 4).Value = "'3.99"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -5.79%  '
$ws.Cells.Item(46, 4).Value = "'2.102.10"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +3.00%  '
$ws.Cells.Item(47, 4).Value = "'0.998"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -0.06%  '
$ws.Cells.Item(48, 4).Value = "'93.68"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +2.52%  '
$ws.Cells.Item(49, 4).Value = "'9.20"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -1.02%  '
$ws.Cells.Item(50, 4).Value = "'108.20"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -2.43%  '
$ws.Cells.Item(51, 4).Value = "'2.847.69"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -0.67%  '
